$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5368310213088989
$ws.Range("B1").Value = 3.569548606872559
$ws.Range("C1").Value = 5.769504070281982
$ws.Range("D1").Value = 1.469729065895081
$ws.Range("E1").Value = 0.8593493103981018
